# Intermediate version of notebook with recall analysis. Partially complete
#
# The results sheet (output_base_ff) gains new metric columns:
#   - "Test Accuracy" is inserted between "Validation Accuracy" and "Model Size"
#   - "Train Recall", "Val Recall", "Drowsy Recall", "Non-Drowsy Recall" are
#     inserted between "Model Size" and "History Plot"
#   - "Train Accuracy" / "Validation Accuracy" values for every run are reset
#     to 0 (placeholder) and the new recall/accuracy metrics are populated.
#
# Net column layout changes from:
#   L=Train Accuracy, M=Validation Accuracy, N=Model Size,
#   O=History Plot, P=Confusion Matrix
# to:
#   L=Train Accuracy, M=Validation Accuracy, N=Test Accuracy, O=Model Size,
#   P=Train Recall, Q=Val Recall, R=Drowsy Recall, S=Non-Drowsy Recall,
#   T=History Plot, U=Confusion Matrix

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# L1 "Train Accuracy" and M1 "Validation Accuracy" are unchanged.
# N1 used to hold "Model Size"; it now holds "Test Accuracy" and the
# "Model Size" header moves one column to the right (O1).
$ws.Range("N1").Value = "Test Accuracy"
$ws.Range("O1").Value = "Model Size"

# P1 used to hold "Confusion Matrix" and O1 used to hold "History Plot";
# both move further right to make room for the four new recall columns.
# Re-use the existing header cell formatting (bold / centered / bordered,
# style index 1, same as every other header cell) by copying it from an
# existing header cell before overwriting the text.
$ws.Range("A1").Copy($ws.Range("Q1")) | Out-Null
$ws.Range("A1").Copy($ws.Range("R1")) | Out-Null
$ws.Range("A1").Copy($ws.Range("S1")) | Out-Null
$ws.Range("A1").Copy($ws.Range("T1")) | Out-Null
$ws.Range("A1").Copy($ws.Range("U1")) | Out-Null

$ws.Range("P1").Value = "Train Recall"
$ws.Range("Q1").Value = "Val Recall"
$ws.Range("R1").Value = "Drowsy Recall"
$ws.Range("S1").Value = "Non-Drowsy Recall"
$ws.Range("T1").Value = "History Plot"
$ws.Range("U1").Value = "Confusion Matrix"

# --- Data rows ----------------------------------------------------------
# For every run (rows 2-9): the old Model Size value shifts from N -> O,
# the old History Plot / Confusion Matrix paths shift from O/P -> T/U, the
# old Train/Validation Accuracy are zeroed out as placeholders, and the
# newly computed Test Accuracy / Train Recall / Val Recall / Drowsy Recall
# / Non-Drowsy Recall values are written in.

# Row 2 (history_1 / confusion_matrix_1)
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 81.99999928474426
$ws.Range("O2").Value = 189.1264686584473
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.96875
$ws.Range("R2").Value = 0.88
$ws.Range("S2").Value = 0.6
$ws.Range("T2").Value = "./plots/base_full_face/hist/history_1.png"
$ws.Range("U2").Value = "./plots/base_full_face/cm/confusion_matrix_1.png"

# Row 3 (history_2 / confusion_matrix_2)
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 80.0000011920929
$ws.Range("O3").Value = 189.1264686584473
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.96875
$ws.Range("R3").Value = 0.84
$ws.Range("S3").Value = 0.62
$ws.Range("T3").Value = "./plots/base_full_face/hist/history_2.png"
$ws.Range("U3").Value = "./plots/base_full_face/cm/confusion_matrix_2.png"

# Row 4 (history_3 / confusion_matrix_3)
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 83.99999737739563
$ws.Range("O4").Value = 126.5555458068848
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.96875
$ws.Range("R4").Value = 0.86
$ws.Range("S4").Value = 0.64
$ws.Range("T4").Value = "./plots/base_full_face/hist/history_3.png"
$ws.Range("U4").Value = "./plots/base_full_face/cm/confusion_matrix_3.png"

# Row 5 (history_4 / confusion_matrix_4)
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 77.99999713897705
$ws.Range("O5").Value = 126.5555458068848
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 0.78
$ws.Range("S5").Value = 0.54
$ws.Range("T5").Value = "./plots/base_full_face/hist/history_4.png"
$ws.Range("U5").Value = "./plots/base_full_face/cm/confusion_matrix_4.png"

# Row 6 (history_5 / confusion_matrix_5)
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 83.99999737739563
$ws.Range("O6").Value = 151.3449745178223
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 0.88
$ws.Range("S6").Value = 0.78
$ws.Range("T6").Value = "./plots/base_full_face/hist/history_5.png"
$ws.Range("U6").Value = "./plots/base_full_face/cm/confusion_matrix_5.png"

# Row 7 (history_6 / confusion_matrix_6)
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 95.99999785423279
$ws.Range("O7").Value = 151.3449745178223
$ws.Range("P7").Value = 1
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 0.98
$ws.Range("S7").Value = 0.6
$ws.Range("T7").Value = "./plots/base_full_face/hist/history_6.png"
$ws.Range("U7").Value = "./plots/base_full_face/cm/confusion_matrix_6.png"

# Row 8 (history_7 / confusion_matrix_7)
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 92.00000166893005
$ws.Range("O8").Value = 101.2740516662598
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 0.9375
$ws.Range("R8").Value = 0.96
$ws.Range("S8").Value = 0.5600000000000001
$ws.Range("T8").Value = "./plots/base_full_face/hist/history_7.png"
$ws.Range("U8").Value = "./plots/base_full_face/cm/confusion_matrix_7.png"

# Row 9 (history_8 / confusion_matrix_8)
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 66.00000262260437
$ws.Range("O9").Value = 101.2740516662598
$ws.Range("P9").Value = 1
$ws.Range("Q9").Value = 0.96875
$ws.Range("R9").Value = 0.7
$ws.Range("S9").Value = 0.5600000000000001
$ws.Range("T9").Value = "./plots/base_full_face/hist/history_8.png"
$ws.Range("U9").Value = "./plots/base_full_face/cm/confusion_matrix_8.png"
